# Add "I0" and "IF" columns (I and J) to the sheet, populating header + 60 data rows (rows 2-61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting from the existing header cell (H1) so the
# new header cells share the same bold/bordered/centered style, then set values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$i0Values = @(9,8,8,7,6,7,7,8,6,6,8,7,8,7,7,8,7,5,5,7,8,7,5,6,4,1,3,8,8,1,6,5,6,6,5,6,6,7,9,8,8,5,6,8,6,7,6,9,10,4,6,2,5,8,6,8,6,8,8,7)
$ifValues = @(9,9,8,8,7,7,8,8,6,7,8,7,8,7,7,8,7,5,5,7,8,7,5,7,4,1,3,8,8,2,6,6,7,6,6,6,7,7,9,9,8,5,6,8,6,7,7,9,10,5,8,2,6,8,7,8,6,8,8,7)

for ($r = 2; $r -le 61; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($r, 10).Value = $ifValues[$idx]
}
